$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying article rows (2-9) were re-sorted into a new order.
# Row 1 (header) and rows 10-11 (unknown/unmatched rows) stay put.

# Note: each row's hyperlink (the rId behind column E) intentionally stays
# anchored to its original cell/target - only the displayed cell text moves
# with the row, mirroring the source diff (no changes to the hyperlink
# relationships were present there).

$ws.Range("A2").Value = "Tunisia attack inquests put back to 2017"
$ws.Range("B2").Value = "2016-03-01T16:31:01UTC"
$ws.Range("C2").Value = 249
$ws.Range("D2").Value = "day_31_beyond"
$ws.Range("E2").Value = "https://www.bbc.com/news/uk-35698927"

$ws.Range("A3").Value = "Tunisia to shut illegal mosques as IS claims deadly hotel attack"
$ws.Range("B3").Value = "2015-06-27T06:07:00UTC"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "day_1"
$ws.Range("E3").Value = "http://www.timesofisrael.com/islamic-state-claims-deadly-hotel-attack-in-tunisia/"

$ws.Range("A4").Value = "PM to chair emergency Cobra meeting after Tunisia attacks"
$ws.Range("B4").Value = "2015-06-27T07:45:35UTC"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "day_1"
$ws.Range("E4").Value = "http://www.itv.com/news/update/2015-06-27/pm-to-chair-emergency-cobra-meeting-after-tunisia-attacks/"

$ws.Range("A5").Value = "Family shocked as Tunis 'break-dance star' becomes mass murderer"
$ws.Range("B5").Value = "2015-06-28T00:13:00UTC"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "day_2_to_30"
$ws.Range("E5").Value = "http://www.timesofisrael.com/tunisia-gunman-was-local-break-dancing-celebrity-fan-of-real-madrid/"

$ws.Range("A6").Value = "How Britain and the EU allowed jihadists to wage war on their own tourists"
$ws.Range("B6").Value = "2015-07-04T21:11:00UTC"
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = "day_2_to_30"
$ws.Range("E6").Value = "https://www.telegraph.co.uk/news/worldnews/africaandindianocean/tunisia/11718577/How-Britain-and-the-EU-allowed-jihadists-to-wage-war-on-their-own-tourists.html"

$ws.Range("A7").Value = "Duke of Sussex unveils Sousse and Bardo Memorial"
$ws.Range("B7").Value = "2019-03-04T12:51:00UTC"
$ws.Range("C7").Value = 1347
$ws.Range("D7").Value = "day_31_beyond"
$ws.Range("E7").Value = "https://www.gov.uk/government/news/duke-of-sussex-unveils-sousse-and-bardo-memorial"

$ws.Range("A8").Value = "Tunisia to close 80 mosques following terror attack"
$ws.Range("B8").Value = "2017-11-07T00:00:00UTC"
$ws.Range("C8").Value = 865
$ws.Range("D8").Value = "day_31_beyond"
$ws.Range("E8").Value = "http://www.dw.com/en/tunisia-to-close-80-mosques-following-terror-attack/a-18544478"

$ws.Range("A9").Value = "Consider Syria IS strikes, defence secretary urges MPs"
$ws.Range("B9").Value = "2015-07-02T23:03:26UTC"
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = "day_2_to_30"
$ws.Range("E9").Value = "https://www.bbc.co.uk/news/uk-33358267"
